$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$excel.ActiveWindow.SmallScroll($null, 7)
